$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J3 was an empty placeholder cell; it now holds the consultant's rating.
$ws.Range("J3").Value = 5

# --- New row 4: another user record ---

# A4 needs the same header/id formatting (border, bold, centered,
# top-aligned) used by A2/A3. Setting .Style directly is a no-op on this
# host, so copy the format from A3 instead.
$ws.Range("A4").Value = 2
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("B4").Value = 3

# C4 and I4 look like a number / a date to Excel's smart-entry parser, so
# force text storage (matching the source data, which stores these as
# plain strings) then drop back to the Normal style so no stray
# number-format style lingers on the cell.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "169707453"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = "youjintyan"
$ws.Range("E4").Value = "YT"
$ws.Range("F4").Value = "None"
$ws.Range("G4").Value = "Арар166"
$ws.Range("H4").Value = "user"

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "2023-11-11"
$ws.Range("I4").Style = "Normal"

# J4 stays empty (consultant_rating not yet set for this user) but the
# cell itself must still exist in the sheet, like J3 did before.
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = ""
$ws.Range("J4").Style = "Normal"
